# TC02_Canine_Filter_Sex-Female.xlsx
# Fixed Tests for SamplePatholoy, SampleType, Se, StageOfDisease, and Study
#
# The "CasesTab" Cypher query in cell B2 of the "startup" sheet is updated
# to drop the trailing `Cohort` column (the OPTIONAL MATCH / co:cohort
# variable is still declared for other reasons, but it is no longer
# returned), matching the rest of the query set which was already
# returning a clean subset of fields.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("startup")

$casesQuery = 'MATCH (s:study)<-[*]-(c:case)<--(demo:demographic)
MATCH (c)<--(diag:diagnosis)
OPTIONAL MATCH (samp:sample)-->(c)
OPTIONAL MATCH (co:cohort)<-[*]-(c)
WITH DISTINCT c, s, demo, diag, co
WHERE demo.sex IN [''Female'']
RETURN  coalesce(c.case_id, '''') AS `Case ID` ,
        coalesce(s.clinical_study_designation, '''') AS `Study Code` ,
        coalesce(s.clinical_study_type, '''') AS  `Study Type`,
        coalesce(demo.breed, '''') AS Breed ,
        coalesce(diag.disease_term, '''') AS Diagnosis ,
        coalesce(diag.stage_of_disease, '''') AS `Stage of Disease` ,
        coalesce(demo.patient_age_at_enrollment, '''') AS Age ,
        coalesce(demo.sex, '''') AS Sex ,
        coalesce(demo.neutered_indicator, '''') AS `Neutered Status`,
        coalesce(demo.weight, '''') AS `Weight (kg)`,
        coalesce(diag.best_response, '''') AS `Response to Treatment`'

$ws.Range("B2").Value = $casesQuery

# The active selection in the saved workbook moved from B4 to B2, and the
# sheet is scrolled back to the top (no frozen/scrolled topLeftCell).
$ws.Range("B2").Select()
